$d = $word.ActiveDocument

# --- 1) "Heureusement pour la star, celle-ci est sur ses gardes apres les
#         dernieres lettres de menace qu'elle a recue." ->
#        "Heureusement pour la star, elle est sur ses gardes apres les
#         lettres de menace qu'elle a recues."
$d.Content.Find.Execute("celle-ci est sur ses gardes", $true, $false, $false, $false, $false, $true, 1, $false, "elle est sur ses gardes", 2) | Out-Null
$d.Content.Find.Execute("les dernières lettres de menace", $true, $false, $false, $false, $false, $true, 1, $false, "les lettres de menace", 2) | Out-Null
$d.Content.Find.Execute("qu’elle a reçue.", $true, $false, $false, $false, $false, $true, 1, $false, "qu’elle a reçues.", 2) | Out-Null

# --- 2) "Pliez les papiers en deux, mélangez-les puis chacun en prend un
#         au hasard." -> "..., mélangez-les, puis chacun en prend un au
#         hasard."
$d.Content.Find.Execute("mélangez-les puis chacun", $true, $false, $false, $false, $false, $true, 1, $false, "mélangez-les, puis chacun", 2) | Out-Null

# --- 3) "... que vous allez jouer mais le stalker se cache toujours !" ->
#        "... que vous allez jouer, mais le stalker se cache toujours !"
$d.Content.Find.Execute("allez jouer mais le stalker", $true, $false, $false, $false, $false, $true, 1, $false, "allez jouer, mais le stalker", 2) | Out-Null

# --- 4) "...etc…)." -> "...etc.)."
$d.Content.Find.Execute("etc…).", $true, $false, $false, $false, $false, $true, 1, $false, "etc.).", 2) | Out-Null

# --- 5) Move the _GoBack bookmark from the end of the "Un fan qui n'est
#         pas interrogé..." paragraph to the end of the "A la fin du tour
#         la star peut accuser..." paragraph, and change that paragraph's
#         closing " :" into "."
$bm = $d.Bookmarks("_GoBack")
$bm.Delete() | Out-Null

$d.Content.Find.Execute("en disant « A l’aide ! Il va me tuer ! » :", $true, $false, $false, $false, $false, $true, 1, $false, "en disant « A l’aide ! Il va me tuer ! ».", 2) | Out-Null

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*accuser un de*") {
        $target = $p
        break
    }
}

# Placing a collapsed bookmark exactly one character before a paragraph's
# end hits an engine edge case, so pad with a throw-away character, anchor
# the bookmark two characters from the (padded) end, then remove the
# padding again. The bookmark stays anchored to the real text.
$r = $target.Range
$r.InsertAfter("X")
$endPos = $target.Range.End

$bmRange = $d.Range($endPos - 2, $endPos - 2)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$padRange = $d.Range($endPos - 2, $endPos - 1)
$padRange.Delete() | Out-Null
